$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("M6").Value = "Jessica S. Tisch"   # Police Commissioner name
$ws.Range("A8").Value = "Volume 31   Number  48"   # Volume/Number line ("47" -> "48")
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"   # week-covering dates

# --- Cells whose underlying style must flip from the "0"/"***.*" placeholder style
#     (s=13) to the ordinary numeric style (s=14), in addition to taking new numbers ---
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)   # paste formats only
$ws.Range("C18").Value = 5
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)   # paste formats only
$ws.Range("C28").Value = 1

# --- Cells that flip the other way: numeric style (s=14/15) -> placeholder text style (s=13) ---
$ws.Range("D23").NumberFormat = "@"   # force text entry so "0" is not coerced to a number
$ws.Range("D23").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D23").PasteSpecial(-4122)   # paste formats only (restores s=13 exactly)
$ws.Range("E23").NumberFormat = "@"   # force text entry so "***.*" is not coerced to a number
$ws.Range("E23").Value = "***.*"
$ws.Range("N23").Copy()
$ws.Range("E23").PasteSpecial(-4122)   # paste formats only (restores s=13 exactly)
$ws.Range("D31").NumberFormat = "@"   # force text entry so "0" is not coerced to a number
$ws.Range("D31").Value = "0"
$ws.Range("C31").Copy()
$ws.Range("D31").PasteSpecial(-4122)   # paste formats only (restores s=13 exactly)
$ws.Range("E31").NumberFormat = "@"   # force text entry so "***.*" is not coerced to a number
$ws.Range("E31").Value = "***.*"
$ws.Range("N31").Copy()
$ws.Range("E31").PasteSpecial(-4122)   # paste formats only (restores s=13 exactly)

# --- Plain numeric value updates across the crime-stats table (rows 15-28) ---
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = -25
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -46.666666666666
$ws.Range("I16").Value = 147
$ws.Range("J16").Value = 161
$ws.Range("K16").Value = -8.695652173913
$ws.Range("L16").Value = -20.967741935483
$ws.Range("M16").Value = -34.080717488789
$ws.Range("N16").Value = 51.546391752577
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 6.666666666666
$ws.Range("I17").Value = 218
$ws.Range("J17").Value = 231
$ws.Range("K17").Value = -5.627705627705
$ws.Range("L17").Value = -13.147410358565
$ws.Range("M17").Value = 29.761904761904
$ws.Range("N17").Value = 194.594594594595
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 106
$ws.Range("J18").Value = 107
$ws.Range("K18").Value = -0.934579439252
$ws.Range("L18").Value = 7.070707070707
$ws.Range("M18").Value = 7.070707070707
$ws.Range("N18").Value = -6.194690265486
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 21.739130434782
$ws.Range("I19").Value = 322
$ws.Range("J19").Value = 333
$ws.Range("K19").Value = -3.303303303303
$ws.Range("L19").Value = 17.090909090909
$ws.Range("M19").Value = 29.838709677419
$ws.Range("N19").Value = 445.762711864407
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("I20").Value = 101
$ws.Range("J20").Value = 133
$ws.Range("K20").Value = -24.060150375939
$ws.Range("L20").Value = -20.472440944881
$ws.Range("M20").Value = 60.31746031746
$ws.Range("N20").Value = 12.222222222222
$ws.Range("C21").Value = 19
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 65
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = -5.797101449275
$ws.Range("I21").Value = 909
$ws.Range("J21").Value = 979
$ws.Range("K21").Value = -7.150153217568
$ws.Range("L21").Value = -4.616998950682
$ws.Range("M21").Value = 10.583941605839
$ws.Range("N21").Value = 104.269662921348
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 9
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -48.275862068965
$ws.Range("L22").Value = -44.444444444444
$ws.Range("M23").Value = 125
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 43.75
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = 60.416666666666
$ws.Range("I24").Value = 856
$ws.Range("J24").Value = 847
$ws.Range("K24").Value = 1.062573789846
$ws.Range("L24").Value = -29.080364540182
$ws.Range("M24").Value = 56.489945155393
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = 133.333333333333
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 190.909090909091
$ws.Range("I25").Value = 286
$ws.Range("J25").Value = 194
$ws.Range("K25").Value = 47.422680412371
$ws.Range("L25").Value = -55.3125
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 300
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = 57.894736842105
$ws.Range("I26").Value = 321
$ws.Range("J26").Value = 342
$ws.Range("K26").Value = -6.140350877192
$ws.Range("L26").Value = -13.709677419354
$ws.Range("M26").Value = -28.507795100222
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 36
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = -10
$ws.Range("L28").Value = -23.404255319148
